# userRegistraion updated: using JavaScript
#
# Adds a new "WMEnvt" column (K) to the envtData sheet, mirroring the
# environment value already stored in column G ("clockServerTarget" /
# "test18"), and updates the sheet view so the new cell is selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in K1 and its value in K2 (same environment tag as G2: "test18")
$ws.Range("K1").Value = "WMEnvt"
$ws.Range("K2").Value = "test18"

# Reflect the new active cell / view position (scrolled so column E is
# leftmost, with K2 selected) just like the edited workbook.
$ws.Range("K2").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 5
